# Apply the ifoCAST error-series update for T55
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing B-column values that changed
$ws.Range("B8").Value = -1.103340811
$ws.Range("B9").Value = -3.010055891
$ws.Range("B13").Value = -0.227799741
$ws.Range("B16").Value = 0.289898428
$ws.Range("B18").Value = 0.08153150799999997
$ws.Range("B19").Value = 0.754884929
$ws.Range("B20").Value = 0.417971554
$ws.Range("B21").Value = 0.691959003
$ws.Range("B22").Value = 0.238543425
$ws.Range("B23").Value = -0.456828469

# Append two new quarters (rows 24 and 25)
$ws.Range("A24").Value = "2025-07-01_diff"
$ws.Range("B24").Value = 0.608342578
$ws.Range("A25").Value = "2025-10-01_diff"

# Match the header-cell style used by the other date labels in column A
$ws.Range("A23").Copy() | Out-Null
$ws.Range("A24:A25").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
